# Update the "Förändrad" (column C) date for every data row on the
# active worksheet from 2023-09-21 (serial 45190) to 2023-09-23
# (serial 45192).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 116) { $lastRow = 116 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value = 45192
    }
}
